$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): extend with two new columns, P and Q, continuing
# the sequential numbering (14, 15). Style them like the rest of the header
# row (bold font, thin box border, centered horizontally, top-aligned
# vertically) so they end up sharing the same cell style as B1:O1.
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

$hdr = $ws.Range("P1:Q1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# --- Data rows (2-25): columns I and K swap (1<->2), columns M and O swap
# (1<->2), and two new columns P and Q are appended with value 2.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2    # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P = 2 (new)
    $ws.Cells.Item($r, 17).Value = 2   # Q = 2 (new)
}
